$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the metric description for the "Taste [6]" tool (row 7, column B):
# replace the old PTMI/PET description with the new TMC/PET description
# (adds the "Average Textual similarity ... (TMC)" metric, i.e. the
# cross-project results row update referenced by the commit message).
$ws.Range("B7").Value = "Average Textual similarity between all the pairs methods called by the test method (TMC)`nProbability of a Class to be affected by Eager Test baisng on its textual content (PET)`n"

# Move the active selection from E8 to B7 and drop the old C1 scroll anchor.
$ws.Range("B7").Select() | Out-Null
